# This script reproduces the commit "Fruta / hortaliza, semanal" for the
# "Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Espinaca" sheet.
#
# The edit inserts one new data row right before the current row 680,
# shifting the former rows 680-731 down to 681-732 (dimension grows from
# A1:R731 to A1:R732), and fills the newly inserted row with a new
# observation (Fecha 45021, Región Metropolitana, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 680; everything from the old row 680 onward
# moves down by one row.
$ws.Rows.Item(680).Insert()

# Fill the newly inserted row 680 with the new record's values.
$ws.Cells.Item(680, 1).Value  = 6
$ws.Cells.Item(680, 2).Value  = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(680, 3).Value  = 'Metropolitana'
$ws.Cells.Item(680, 4).Value  = 45021
$ws.Cells.Item(680, 5).Value  = 13
$ws.Cells.Item(680, 6).Value  = 100112012
$ws.Cells.Item(680, 7).Value  = 'Espinaca'
$ws.Cells.Item(680, 8).Value  = 'Sin especificar'
$ws.Cells.Item(680, 9).Value  = 'Primera'
$ws.Cells.Item(680, 10).Value = 450
$ws.Cells.Item(680, 11).Value = 8000
$ws.Cells.Item(680, 12).Value = 9000
$ws.Cells.Item(680, 13).Value = 8578
$ws.Cells.Item(680, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(680, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(680, 16).Value = 858
$ws.Cells.Item(680, 17).Value = 10
$ws.Cells.Item(680, 18).Value = 'Hortaliza'
